$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2,1,2575),
    @(3,1,2594),
    @(4,1,2595),
    @(5,1,2599),
    @(6,1,2601),
    @(7,1,2604),
    @(8,1,2607),
    @(9,1,2608),
    @(10,1,2610),
    @(11,1,2612),
    @(12,1,2613),
    @(13,1,2616),
    @(14,1,2617),
    @(15,1,2646),
    @(16,2,2316),
    @(17,2,2317),
    @(18,2,2323),
    @(19,2,2341),
    @(20,2,2348),
    @(21,2,2366),
    @(22,2,2406),
    @(23,2,2536),
    @(24,2,2680),
    @(25,3,2449),
    @(26,3,2694),
    @(27,3,2732),
    @(28,3,2738),
    @(29,3,3015),
    @(30,3,3058),
    @(31,3,3082),
    @(32,5,2792),
    @(33,5,2794),
    @(34,5,2796),
    @(35,5,2807),
    @(36,5,2835),
    @(37,5,2843),
    @(38,5,2847),
    @(39,5,2849),
    @(40,5,2854),
    @(41,5,2860),
    @(42,5,2861),
    @(43,5,2862),
    @(44,5,2874),
    @(45,5,2885),
    @(46,5,2890),
    @(47,5,2904),
    @(48,5,2905),
    @(49,5,2933),
    @(50,5,2945),
    @(51,5,2958),
    @(52,5,2960),
    @(53,5,2998),
    @(54,6,2651),
    @(55,6,2652),
    @(56,6,2974),
    @(57,6,2979),
    @(58,6,2980),
    @(59,6,3046),
    @(60,6,3083),
    @(61,6,3095),
    @(62,7,2817),
    @(63,7,2824),
    @(64,7,2833),
    @(65,7,2837),
    @(66,8,2545),
    @(67,8,2549),
    @(68,8,2591),
    @(69,8,2593),
    @(70,8,2605),
    @(71,8,2639),
    @(72,8,2775),
    @(73,8,2777),
    @(74,8,2780),
    @(75,8,2783),
    @(76,8,2786),
    @(77,8,2845),
    @(78,9,2285),
    @(79,9,2292),
    @(80,10,2464),
    @(81,10,2466),
    @(82,10,2475),
    @(83,10,2514),
    @(84,11,2484),
    @(85,11,2488),
    @(86,11,2492),
    @(87,11,2516),
    @(88,12,2397),
    @(89,12,2461),
    @(90,13,2625),
    @(91,13,2626),
    @(92,13,2720),
    @(93,13,2721),
    @(94,13,2729),
    @(95,14,1865),
    @(96,14,2084),
    @(97,14,2086),
    @(98,14,2198),
    @(99,14,2239),
    @(100,15,1011),
    @(101,16,5596),
    @(102,17,3176),
    @(103,17,3177),
    @(104,17,3178),
    @(105,17,3179),
    @(106,17,3191),
    @(107,17,3194),
    @(108,17,3198),
    @(109,17,3204),
    @(110,18,5466),
    @(111,18,5472),
    @(112,18,5475),
    @(113,18,5547),
    @(114,18,5554),
    @(115,19,501),
    @(116,20,5148),
    @(117,21,5320),
    @(118,21,5323),
    @(119,22,221),
    @(120,22,313),
    @(121,23,4069),
    @(122,23,4320),
    @(123,24,1675),
    @(124,25,3233),
    @(125,25,3297),
    @(126,26,1262),
    @(127,27,5051),
    @(128,28,4585),
    @(129,29,3463),
    @(130,29,3525),
    @(131,29,3593),
    @(132,29,3618),
    @(133,29,3626),
    @(134,29,3627),
    @(135,29,3640),
    @(136,29,3912),
    @(137,29,3927),
    @(138,29,3929),
    @(139,30,372)
)

foreach ($row in $rows) {
    $r = $row[0]
    $a = $row[1]
    $b = $row[2]
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
}

$ws.Range("A87:B139").NumberFormat = "0"
